# Implementation of initial AJAX handlers and dialogs
# Populate the Message / Level columns for existing rows 10008, 10010, 10011
# and append two new rows for 10022 and 10023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 24 -> Code 10022 (written first so the shared string table keeps
# the same ordering as produced by the original authoring tool)
$ws.Range("A24").Value = 10022
$ws.Range("B24").Value = "message_10022_person_record_deleted_successfully"
$ws.Range("D24").Value = "Success"

# New row 25 -> Code 10023
$ws.Range("A25").Value = 10023
$ws.Range("B25").Value = "message_10023_person_multiple_records_deleted_successfully"
$ws.Range("D25").Value = "Success"

# Row 10 -> Code 10008
$ws.Range("B10").Value = "message_10008_ajax_login_successful"
$ws.Range("D10").Value = "Success"

# Row 12 -> Code 10010
$ws.Range("B12").Value = "message_10010_invalid_datetime_format"
$ws.Range("D12").Value = "Error"

# Row 13 -> Code 10011
$ws.Range("B13").Value = "message_10011_expired_session_renewed"
$ws.Range("D13").Value = "Information"

# Update the selected cell to match the saved view state
$ws.Range("C16").Select()
